$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The working-experience sheet lists bullet points ("why" column, E) that each
# ended with a trailing period. This commit removes those trailing periods.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 19 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value()
    if ($null -ne $val -and $val -is [string] -and $val.EndsWith(".")) {
        $cell.Value = $val.Substring(0, $val.Length - 1)
    }
}

$ws.Range("E24").Select()
